# Append new attendance/"Entrada PM" records to the "Actual" worksheet.
# The source data was produced programmatically (inline strings), so every
# value here -- including the numeric-looking IDs and the date/time text --
# must be written as literal TEXT, not auto-converted to a number/date by
# Excel's input parser. We force that by pre-formatting the target range as
# Text ("@") before assigning values, then clearing the formatting again so
# the cells end up with the workbook's default (unstyled) look, matching the
# rest of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Actual")

$newRows = @(
    @{ A = "1053868254"; B = "Tatiana Pachon"; C = "Administrativa"; D = "2024-02-07"; E = "13:28:28"; F = "13:28:28"; G = "13:28:28"; H = "Entrada PM" },
    @{ A = "1015463003"; B = "Leonardo Maje";  C = "Administrativa"; D = "2024-02-07"; E = "13:36:18"; F = "13:36:18"; G = "13:36:18"; H = "Entrada PM" },
    @{ A = "1054398414"; B = "Julian Largo";   C = "Administrativa"; D = "2024-02-07"; E = "13:41:18"; F = "13:41:18"; G = "13:41:18"; H = "Entrada PM" },
    @{ A = "10101010";   B = "Proveedor";      C = "Tercero";        D = "2024-02-07"; E = "14:13:10"; F = "14:13:10"; G = "14:13:10"; H = "Entrada PM" },
    @{ A = "1054398414"; B = "Julian Largo";   C = "Administrativa"; D = "2024-02-07"; E = "14:15:28"; F = "13:41:18"; G = "14:15:28"; H = "Entrada PM" }
)

$startRow = 4
$endRow = $startRow + $newRows.Count - 1
$fullRange = $ws.Range("A" + $startRow + ":H" + $endRow)
$fullRange.NumberFormat = "@"

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Range("A$r").Value = $row.A
    $ws.Range("B$r").Value = $row.B
    $ws.Range("C$r").Value = $row.C
    $ws.Range("D$r").Value = $row.D
    $ws.Range("E$r").Value = $row.E
    $ws.Range("F$r").Value = $row.F
    $ws.Range("G$r").Value = $row.G
    $ws.Range("H$r").Value = $row.H
}

$fullRange.ClearFormats()
